$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1505.730094935563
$ws.Range("C2").Value = 1482.659367292297
$ws.Range("D2").Value = 1054.218203400688
$ws.Range("E2").Value = 1641.065310037129
$ws.Range("F2").Value = 1540.366991430833
$ws.Range("G2").Value = 1589.319899554724
$ws.Range("H2").Value = 1642.72015540434
$ws.Range("B3").Value = 1543.621119192856
$ws.Range("C3").Value = 1489.660156732298
$ws.Range("D3").Value = 1028.556352226926
$ws.Range("E3").Value = 1649.511916152012
$ws.Range("F3").Value = 1565.302392880995
$ws.Range("G3").Value = 1593.901005905419
$ws.Range("H3").Value = 1650.71200068266
$ws.Range("B4").Value = 1421.612905331771
$ws.Range("C4").Value = 1536.745462111212
$ws.Range("D4").Value = 770.0498650957028
$ws.Range("E4").Value = 1628.84884246994
$ws.Range("F4").Value = 1457.830666557993
$ws.Range("G4").Value = 1555.648480433783
$ws.Range("H4").Value = 1629.731575774725
$ws.Range("B5").Value = 1599.557046887749
$ws.Range("C5").Value = 1589.172341735097
$ws.Range("D5").Value = 372.791991144114
$ws.Range("E5").Value = 1665.468823484424
$ws.Range("F5").Value = 1601.191988968261
$ws.Range("G5").Value = 1604.157106965235
$ws.Range("H5").Value = 1665.656478443706
$ws.Range("B6").Value = 1585.162946869089
$ws.Range("C6").Value = 1586.12702903763
$ws.Range("D6").Value = 187.6397414329435
$ws.Range("E6").Value = 1662.538223456871
$ws.Range("F6").Value = 1587.64886275616
$ws.Range("G6").Value = 1592.155716468016
$ws.Range("H6").Value = 1662.592171833827
$ws.Range("B7").Value = 1612.1304178267
$ws.Range("C7").Value = 1595.448448823696
$ws.Range("D7").Value = 341.2147995783006
$ws.Range("E7").Value = 1666.496209868257
$ws.Range("F7").Value = 1614.422536682491
$ws.Range("G7").Value = 1605.734104783698
$ws.Range("H7").Value = 1666.659364228769
$ws.Range("B8").Value = 1524.199078212293
$ws.Range("C8").Value = 1433.203986887945
$ws.Range("D8").Value = 326.4477621688607
$ws.Range("E8").Value = 1621.016637854812
$ws.Range("F8").Value = 1524.716381781736
$ws.Range("G8").Value = 1465.874424671512
$ws.Range("H8").Value = 1621.249539289433
$ws.Range("B9").Value = 1625.195807347984
$ws.Range("C9").Value = 1546.589363101275
$ws.Range("D9").Value = 427.8133966681748
$ws.Range("E9").Value = 1663.019301873826
$ws.Range("F9").Value = 1628.197875621464
$ws.Range("G9").Value = 1566.306438240137
$ws.Range("H9").Value = 1663.386026970141
$ws.Range("B10").Value = 1570.125201648339
$ws.Range("C10").Value = 1589.572658536279
$ws.Range("D10").Value = 1009.221651457342
$ws.Range("E10").Value = 1662.741626091052
$ws.Range("F10").Value = 1596.381014948465
$ws.Range("G10").Value = 1641.04387868246
$ws.Range("H10").Value = 1664.172647136288
$ws.Range("B11").Value = 1555.234797529095
$ws.Range("C11").Value = 1591.98943187592
$ws.Range("D11").Value = 916.8538657187961
$ws.Range("E11").Value = 1664.139707617842
$ws.Range("F11").Value = 1580.998156961288
$ws.Range("G11").Value = 1637.111851400533
$ws.Range("H11").Value = 1665.265662253962
$ws.Range("B12").Value = 1297.715867080148
$ws.Range("C12").Value = 1594.044308298733
$ws.Range("D12").Value = 477.0849773684319
$ws.Range("E12").Value = 1650.286811872499
$ws.Range("F12").Value = 1344.132768786582
$ws.Range("G12").Value = 1605.589635760548
$ws.Range("H12").Value = 1650.444086449992
$ws.Range("B13").Value = 1560.211442101217
$ws.Range("C13").Value = 1587.962461674932
$ws.Range("D13").Value = 881.0403330087587
$ws.Range("E13").Value = 1660.983155509591
$ws.Range("F13").Value = 1583.751461116511
$ws.Range("G13").Value = 1631.387590119088
$ws.Range("H13").Value = 1662.183841794597
